$d = $word.ActiveDocument

# Locate the "Objetivos" heading paragraph (List Paragraph, numbered) that
# precedes the (currently empty) paragraph before "Escopo".
$objetivosIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Objetivos") {
        $objetivosIndex = $i
    }
    $i = $i + 1
}

$objPara = $d.Paragraphs.Item($objetivosIndex)

# Insert a brand-new paragraph right after the "Objetivos" heading.
$objPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($objetivosIndex + 1)

# The new paragraph is a regular body paragraph (no list/heading style),
# indented and justified, matching the other narrative paragraphs in the doc.
$newPara.Range.Style = "Normal"
$newPara.Range.ListFormat.ConvertNumbersToText()
$newPara.LeftIndent = 18
$newPara.Alignment = 3

$bodyText = "A empresa parceira possui um aplicativo que realiza estatísticas sobre o plantio e colheita de plantas de soja. Este aplicativo é bastante acessado por agricultores brasileiros e eles sentem falta de uma inovação no aplicativo. Sendo assim, o objetivo é realizar a implementação de uma nova interface para o envio de amostras de algumas plantas de soja da plantação, de forma que a contagem de vagens seja feita de forma automática e o agricultor não necessite contá-las (o que é feito hoje). O projeto terá início no dia 29/08/2022. O orçamento máximo para cada integrante, por sprint, deverá ser no máximo de R`$1.071,00. Para cada sprint, o valor disponibilizado será de R`$7.500,00. O projeto será considerado um sucesso se atender a todos os critérios de aceitação das entregas, respeitar as restrições e cumprir o cronograma de execução, com data de entrega final datada em 29/11/2022."

$newPara.Range.Text = $bodyText
